# "merged chapters 2, 3 and 4"
#
# 1. Drop the leftover "  Types of algorithms(basic notions)" strike-through
#    note that trailed the "Applications in biology and medicine" bullet.
# 2. Recolor the "The role of CASP" / "Types of protein structure
#    prediction" / "Template free modeling" / "Template based modeling"
#    bullets red, matching the rest of the (now-merged) chapter.
# 3. Retitle "Methods for determining tertiary structure" to
#    "Methods for determining protein structure".
# 4. Remove the now-redundant standalone "Methods for determining
#    secondary structure" bullet.

$d = $word.ActiveDocument

# --- 1: trim the paragraph "Applications in biology and medicine  Types of
# algorithms(basic notions)" back down to "...and medicine".
$p = $d.Paragraphs(15)
$full = $p.Range.Text
$idx = $full.IndexOf("  Types of algorithms")
$tailStart = $p.Range.Start + $idx
$tailRange = $d.Range($tailStart, $p.Range.End - 1)
$tailRange.Delete()

# --- 2: color the four "protein structure prediction" intro bullets red.
foreach ($i in 19, 20, 21, 22) {
    $d.Paragraphs($i).Range.Font.Color = 255
}

# --- 3: "tertiary" -> "protein" in "Methods for determining tertiary
# structure", splitting the run (nudge+revert Bold forces the run break)
# so the final XML has three runs just like the rest of the authored doc.
$p = $d.Paragraphs(23)
$full = $p.Range.Text
$wordStart = $p.Range.Start + $full.IndexOf("tertiary")
$wordEnd = $wordStart + [string]"tertiary".Length
$wordRange = $d.Range($wordStart, $wordEnd)
$wordRange.Text = "protein"
$newEnd = $wordStart + [string]"protein".Length
$newRange = $d.Range($wordStart, $newEnd)
$newRange.Bold = 1
$newRange.Bold = 0

# --- 4: delete the standalone "Methods for determining secondary
# structure" bullet entirely (folded into the renamed bullet above).
$d.Paragraphs(32).Range.Delete()
